# Apply scheduled-runner price/profit updates to the Coeurl_Profits workbook.
# Each sheet corresponds to a crafting/gathering job; columns H-N hold the
# currentAveragePrice*, LevePrice*, and LeveProfit* computed figures that were
# refreshed by the runner. Cells that had no prior value (and are left blank
# by the refresh) are cleared so they do not linger with stale numbers.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = ""
$ws.Range("H43").Value = 23664.924
$ws.Range("I43").Value = 32368.334
$ws.Range("J43").Value = 16204.857
$ws.Range("K43").Value = 32368.334
$ws.Range("L43").Value = 16204.857
$ws.Range("M43").Value = -32299.334
$ws.Range("N43").Value = -16342.857
$ws.Range("H70").Value = 3463.4546
$ws.Range("I70").Value = 2550
$ws.Range("K70").Value = 7650
$ws.Range("M70").Value = -7380
$ws.Range("H73").Value = 3463.4546
$ws.Range("I73").Value = 2550
$ws.Range("K73").Value = 7650
$ws.Range("M73").Value = -6714
$ws.Range("H86").Value = 4064.1428
$ws.Range("I86").Value = 3727.0908
$ws.Range("K86").Value = 3727.0908
$ws.Range("M86").Value = -2604.0908
$ws.Range("H89").Value = 4064.1428
$ws.Range("I89").Value = 3727.0908
$ws.Range("K89").Value = 18635.454
$ws.Range("M89").Value = -13019.454
$ws.Range("H100").Value = 2291.4443
$ws.Range("I100").Value = 1946.2858
$ws.Range("J100").Value = 3499.5
$ws.Range("K100").Value = 1946.2858
$ws.Range("L100").Value = 3499.5
$ws.Range("M100").Value = -1405.2858
$ws.Range("N100").Value = -4581.5
$ws.Range("H111").Value = 3000
$ws.Range("I111").Value = 3000
$ws.Range("K111").Value = 9000
$ws.Range("M111").Value = -5933
$ws.Range("H141").Value = 1116.5
$ws.Range("I141").Value = 1116.5
$ws.Range("K141").Value = 3349.5
$ws.Range("M141").Value = 1830.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3570.1956
$ws.Range("I32").Value = 2734.9736
$ws.Range("K32").Value = 2734.9736
$ws.Range("M32").Value = -2447.9736
$ws.Range("H43").Value = 38492.2
$ws.Range("J43").Value = 38115.25
$ws.Range("L43").Value = 38115.25
$ws.Range("N43").Value = -38741.25
$ws.Range("H117").Value = 40081.332
$ws.Range("J117").Value = 40081.332
$ws.Range("L117").Value = 40081.332
$ws.Range("N117").Value = -49259.332
$ws.Range("H122").Value = 85216.5
$ws.Range("I122").Value = 111853.78
$ws.Range("K122").Value = 335561.34
$ws.Range("M122").Value = -333111.34
$ws.Range("H132").Value = 2978
$ws.Range("I132").Value = 2064
$ws.Range("J132").Value = 5720
$ws.Range("K132").Value = 6192
$ws.Range("L132").Value = 17160
$ws.Range("M132").Value = -3662
$ws.Range("N132").Value = -22220

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2129.5454
$ws.Range("I134").Value = 1564.4231
$ws.Range("K134").Value = 4693.2693
$ws.Range("M134").Value = -2158.2693

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 121211
$ws.Range("J31").Value = 13483.333
$ws.Range("L31").Value = 13483.333
$ws.Range("N31").Value = -14073.333
$ws.Range("H34").Value = 121211
$ws.Range("J34").Value = 13483.333
$ws.Range("L34").Value = 13483.333
$ws.Range("N34").Value = -13887.333
$ws.Range("H62").Value = 6368.1
$ws.Range("I62").Value = 6996
$ws.Range("J62").Value = 5426.25
$ws.Range("K62").Value = 6996
$ws.Range("L62").Value = 5426.25
$ws.Range("M62").Value = -6372
$ws.Range("N62").Value = -6674.25
$ws.Range("H65").Value = 6368.1
$ws.Range("I65").Value = 6996
$ws.Range("J65").Value = 5426.25
$ws.Range("K65").Value = 34980
$ws.Range("L65").Value = 27131.25
$ws.Range("M65").Value = -31860
$ws.Range("N65").Value = -33371.25
$ws.Range("H99").Value = 2747.5
$ws.Range("I99").Value = 2276.3635
$ws.Range("J99").Value = 4475
$ws.Range("K99").Value = 2276.3635
$ws.Range("L99").Value = 4475
$ws.Range("M99").Value = -778.3634999999999
$ws.Range("N99").Value = -7471
$ws.Range("H107").Value = 949.8095
$ws.Range("I107").Value = 958.5714
$ws.Range("K107").Value = 958.5714
$ws.Range("M107").Value = 961.4286
$ws.Range("H126").Value = 2747.5
$ws.Range("I126").Value = 2276.3635
$ws.Range("J126").Value = 4475
$ws.Range("K126").Value = 6829.0905
$ws.Range("L126").Value = 13425
$ws.Range("M126").Value = -4359.0905
$ws.Range("N126").Value = -18365

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 5865.1665
$ws.Range("J104").Value = 5865.1665
$ws.Range("L104").Value = 17595.4995
$ws.Range("N104").Value = -22837.4995
$ws.Range("H136").Value = 1007812.8
$ws.Range("I136").Value = 2001025.6
$ws.Range("K136").Value = 6003076.800000001
$ws.Range("M136").Value = -5997976.800000001
$ws.Range("H137").Value = 3095.818
$ws.Range("I137").Value = 2936.4285
$ws.Range("J137").Value = 3374.75
$ws.Range("K137").Value = 8809.2855
$ws.Range("L137").Value = 10124.25
$ws.Range("M137").Value = -3709.2855
$ws.Range("N137").Value = -20324.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6137.5557
$ws.Range("I80").Value = 3402.5
$ws.Range("J80").Value = 6919
$ws.Range("K80").Value = 3402.5
$ws.Range("L80").Value = 6919
$ws.Range("M80").Value = -2404.5
$ws.Range("N80").Value = -8915
$ws.Range("H83").Value = 6137.5557
$ws.Range("I83").Value = 3402.5
$ws.Range("J83").Value = 6919
$ws.Range("K83").Value = 17012.5
$ws.Range("L83").Value = 34595
$ws.Range("M83").Value = -12020.5
$ws.Range("N83").Value = -44579
$ws.Range("H107").Value = 1669.3334
$ws.Range("I107").Value = 1462.3334
$ws.Range("J107").Value = 2083.3333
$ws.Range("K107").Value = 1462.3334
$ws.Range("L107").Value = 2083.3333
$ws.Range("M107").Value = 457.6666
$ws.Range("N107").Value = -5923.3333
$ws.Range("H113").Value = 7916.3335
$ws.Range("I113").Value = 4874.5
$ws.Range("K113").Value = 4874.5
$ws.Range("M113").Value = -2704.5
$ws.Range("H122").Value = 3393.8667
$ws.Range("I122").Value = 2917.6667
$ws.Range("J122").Value = 5298.6665
$ws.Range("K122").Value = 8753.000100000001
$ws.Range("L122").Value = 15895.9995
$ws.Range("M122").Value = -6303.000100000001
$ws.Range("N122").Value = -20795.9995
$ws.Range("H126").Value = 36030.7
$ws.Range("I126").Value = 48686.715
$ws.Range("K126").Value = 146060.145
$ws.Range("M126").Value = -143590.145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I55").Value = 297
$ws.Range("J55").Value = 492
$ws.Range("K55").Value = 297
$ws.Range("L55").Value = 492
$ws.Range("M55").Value = -124
$ws.Range("N55").Value = -838
$ws.Range("H82").Value = 2273.4
$ws.Range("I82").Value = 2097.3333
$ws.Range("K82").Value = 2097.3333
$ws.Range("M82").Value = -1736.3333
$ws.Range("H85").Value = 2273.4
$ws.Range("I85").Value = 2097.3333
$ws.Range("K85").Value = 2097.3333
$ws.Range("M85").Value = -849.3332999999998
$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180
$ws.Range("H125").Value = 105624.125
$ws.Range("J125").Value = 495000
$ws.Range("L125").Value = 495000
$ws.Range("N125").Value = -504840
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""
$ws.Range("H132").Value = 4169.6665
$ws.Range("I132").Value = 3251.2666
$ws.Range("J132").Value = 6465.6665
$ws.Range("K132").Value = 9753.799800000001
$ws.Range("L132").Value = 19396.9995
$ws.Range("M132").Value = -7223.799800000001
$ws.Range("N132").Value = -24456.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 353.13333
$ws.Range("J113").Value = 510
$ws.Range("L113").Value = 1530
$ws.Range("N113").Value = -5870
$ws.Range("H126").Value = 3370.9546
$ws.Range("I126").Value = 2034.75
$ws.Range("K126").Value = 6104.25
$ws.Range("M126").Value = -3634.25
$ws.Range("H132").Value = 4695.0605
$ws.Range("I132").Value = 4591.7812
$ws.Range("K132").Value = 13775.3436
$ws.Range("M132").Value = -11245.3436

